# Inserts a new weekly price record as row 171 (Arándano (blue) / O'Higgins),
# pushing the existing rows 171-207 down to 172-208.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 171; everything below shifts down one row.
$ws.Rows(171).Insert()

# Populate the new row 171 with the new record's data.
$ws.Range("A171").Value = 6
$ws.Range("B171").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C171").Value = "Metropolitana"
$ws.Range("D171").Value = 44511
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100101
$ws.Range("H171").Value = "Berries"
$ws.Range("I171").Value = 100101001
$ws.Range("J171").Value = "Arándano (blue)"
$ws.Range("K171").Value = "Sin especificar"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 3500
$ws.Range("N171").Value = 6000
$ws.Range("O171").Value = 6000
$ws.Range("P171").Value = 6000
$ws.Range("Q171").Value = "$/bandeja 2 kilos"
$ws.Range("R171").Value = "Región de O'Higgins"
$ws.Range("S171").Value = 3000
$ws.Range("T171").Value = 2
